$wb = $excel.ActiveWorkbook

# Rename the "Los Angeles" sheet to "Los_Angeles"
$ws = $wb.Worksheets.Item("Los Angeles")
$ws.Name = "Los_Angeles"

# Fix the misspelled latitude headers (column C = Origin, column I = Destination)
$ws.Range("C1").Value = "Origin Latitude"
$ws.Range("I1").Value = "Destination Latitude"

# Update the active selection shown in the sheet to G2
$ws.Activate()
$ws.Range("G2").Select()
